$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.41718995543362
$ws.Range("D2").Value = 7.755181548233864
$ws.Range("E2").Value = 14.49979653917049
$ws.Range("F2").Value = 43.97122994501682
$ws.Range("G2").Value = 52.82079700436275
$ws.Range("H2").Value = 20.2283635944198
$ws.Range("J2").Value = 11.49507952569678
$ws.Range("K2").Value = 11.33238299099355
$ws.Range("L2").Value = 9.564494407665732
$ws.Range("B3").Value = 17.34635167588763
$ws.Range("D3").Value = 7.762401211725274
$ws.Range("E3").Value = 14.52282562918954
$ws.Range("F3").Value = 43.97058517885589
$ws.Range("G3").Value = 52.7298647210302
$ws.Range("H3").Value = 20.2555147304697
$ws.Range("J3").Value = 11.51216733116932
$ws.Range("K3").Value = 11.16613743602768
$ws.Range("L3").Value = 9.548551086356921
$ws.Range("B4").Value = 17.30634870736386
$ws.Range("D4").Value = 7.76745505799439
$ws.Range("E4").Value = 14.53798551284394
$ws.Range("F4").Value = 43.97874373509062
$ws.Range("G4").Value = 52.68614867602369
$ws.Range("H4").Value = 20.27534863231163
$ws.Range("J4").Value = 11.52324755797173
$ws.Range("K4").Value = 11.06571748400558
$ws.Range("L4").Value = 9.54042635671566
$ws.Range("B5").Value = 17.29093849271833
$ws.Range("D5").Value = 7.769671107950447
$ws.Range("E5").Value = 14.54442026516896
$ws.Range("F5").Value = 43.98421924969553
$ws.Range("G5").Value = 52.6713883395638
$ws.Range("H5").Value = 20.28422568351464
$ws.Range("J5").Value = 11.52791113959607
$ws.Range("K5").Value = 11.02526648881118
$ws.Range("L5").Value = 9.537536910430573
$ws.Range("B6").Value = 17.2884338371441
$ws.Range("D6").Value = 7.77004854996454
$ws.Range("E6").Value = 14.54550428657571
$ws.Range("F6").Value = 43.98525833245963
$ws.Range("G6").Value = 52.66912203389522
$ws.Range("H6").Value = 20.28574767583692
$ws.Range("J6").Value = 11.52869449187383
$ws.Range("K6").Value = 11.01857966134541
$ws.Range("L6").Value = 9.537082653156329
$ws.Range("B7").Value = 17.30613725423886
$ws.Range("D7").Value = 7.767484309913996
$ws.Range("E7").Value = 14.53807125299095
$ws.Range("F7").Value = 43.97880887231177
$ws.Range("G7").Value = 52.68593723847675
$ws.Range("H7").Value = 20.27546513530471
$ws.Range("J7").Value = 11.52330985169973
$ws.Range("K7").Value = 11.06516996911823
$ws.Range("L7").Value = 9.540385678612912
$ws.Range("B8").Value = 17.39204839943698
$ws.Range("D8").Value = 7.757542267282795
$ws.Range("E8").Value = 14.50752558811363
$ws.Range("F8").Value = 43.96923229691138
$ws.Range("G8").Value = 52.78693399194715
$ws.Range("H8").Value = 20.23706839118352
$ws.Range("J8").Value = 11.50084954007952
$ws.Range("K8").Value = 11.27475087855665
$ws.Range("L8").Value = 9.558653491098115
$ws.Range("B9").Value = 17.58759902592701
$ws.Range("D9").Value = 7.742954329081633
$ws.Range("E9").Value = 14.45569564103225
$ws.Range("F9").Value = 44.01830183031989
$ws.Range("G9").Value = 53.08068679296596
$ws.Range("H9").Value = 20.18690087028645
$ws.Range("J9").Value = 11.46145502366601
$ws.Range("K9").Value = 11.69645510957863
$ws.Range("L9").Value = 9.607545656719656
$ws.Range("B10").Value = 17.74687098597473
$ws.Range("D10").Value = 7.735204152027969
$ws.Range("E10").Value = 14.42250480894243
$ws.Range("F10").Value = 44.09562472542324
$ws.Range("G10").Value = 53.35408787443846
$ws.Range("H10").Value = 20.16539917269507
$ws.Range("J10").Value = 11.43532233504487
$ws.Range("K10").Value = 12.00958546226404
$ws.Range("L10").Value = 9.651232107402624
$ws.Range("B11").Value = 17.82250321129718
$ws.Range("D11").Value = 7.732317448911892
$ws.Range("E11").Value = 14.40846039829864
$ws.Range("F11").Value = 44.13972086484745
$ws.Range("G11").Value = 53.49074981044335
$ws.Range("H11").Value = 20.15895706942327
$ws.Range("J11").Value = 11.42403899071069
$ws.Range("K11").Value = 12.15208371479825
$ws.Range("L11").Value = 9.672743652304112
$ws.Range("B12").Value = 17.85158108539466
$ws.Range("D12").Value = 7.731315742468878
$ws.Range("E12").Value = 14.403293240997
$ws.Range("F12").Value = 44.15769623033434
$ws.Range("G12").Value = 53.54424451572481
$ws.Range("H12").Value = 20.15699786864088
$ws.Range("J12").Value = 11.41985282652525
$ws.Range("K12").Value = 12.20599764831363
$ws.Range("L12").Value = 9.681120534010372
$ws.Range("B13").Value = 17.8452995007156
$ws.Range("D13").Value = 7.731527418886914
$ws.Range("E13").Value = 14.40439936513353
$ws.Range("F13").Value = 44.15376821822537
$ws.Range("G13").Value = 53.53264632078655
$ws.Range("H13").Value = 20.15739845670542
$ws.Range("J13").Value = 11.42075054599284
$ws.Range("K13").Value = 12.19438926438189
$ws.Range("L13").Value = 9.679306229198787
$ws.Range("B14").Value = 17.8248868096473
$ws.Range("D14").Value = 7.732233208198435
$ws.Range("E14").Value = 14.40803226630482
$ws.Range("F14").Value = 44.14117415098956
$ws.Range("G14").Value = 53.49511604797748
$ws.Range("H14").Value = 20.15878625910571
$ws.Range("J14").Value = 11.4236928589367
$ws.Range("K14").Value = 12.15652049719539
$ws.Range("L14").Value = 9.673428222946011
$ws.Range("B15").Value = 17.81243984381475
$ws.Range("D15").Value = 7.732677417455967
$ws.Range("E15").Value = 14.41027719498864
$ws.Range("F15").Value = 44.13362604754797
$ws.Range("G15").Value = 53.47235402222095
$ws.Range("H15").Value = 20.15969887597218
$ws.Range("J15").Value = 11.4255063775665
$ws.Range("K15").Value = 12.13331705937215
$ws.Range("L15").Value = 9.669857704894657
$ws.Range("B16").Value = 17.74199055453464
$ws.Range("D16").Value = 7.735405622979481
$ws.Range("E16").Value = 14.42344382054138
$ws.Range("F16").Value = 44.09292199452135
$ws.Range("G16").Value = 53.34540205710351
$ws.Range("H16").Value = 20.16588738534665
$ws.Range("J16").Value = 11.43607186405414
$ws.Range("K16").Value = 12.00026922413401
$ws.Range("L16").Value = 9.649858852291075
$ws.Range("B17").Value = 17.69957268274808
$ws.Range("D17").Value = 7.737242608487964
$ws.Range("E17").Value = 14.43179082330875
$ws.Range("F17").Value = 44.07023264496946
$ws.Range("G17").Value = 53.27065309262977
$ws.Range("H17").Value = 20.17053920648361
$ws.Range("J17").Value = 11.44270804407102
$ws.Range("K17").Value = 11.91862291967657
$ws.Range("L17").Value = 9.638006589423943
$ws.Range("B18").Value = 17.67547544576712
$ws.Range("D18").Value = 7.738359356582304
$ws.Range("E18").Value = 14.4366910574779
$ws.Range("F18").Value = 44.05802249096931
$ws.Range("G18").Value = 53.22881787655827
$ws.Range("H18").Value = 20.17352910871099
$ws.Range("J18").Value = 11.44658192221359
$ws.Range("K18").Value = 11.87167035526705
$ws.Range("L18").Value = 9.631343987659708
$ws.Range("B19").Value = 17.66736871141697
$ws.Range("D19").Value = 7.738747814770064
$ws.Range("E19").Value = 14.438367254015
$ws.Range("F19").Value = 44.05403279872526
$ws.Range("G19").Value = 53.21485282385969
$ws.Range("H19").Value = 20.17459541272354
$ws.Range("J19").Value = 11.44790333775678
$ws.Range("K19").Value = 11.85577613462062
$ws.Range("L19").Value = 9.629114819733886
$ws.Range("B20").Value = 17.70405718611956
$ws.Range("D20").Value = 7.737040834864853
$ws.Range("E20").Value = 14.43089200093395
$ws.Range("F20").Value = 44.07256104985124
$ws.Range("G20").Value = 53.27849051843191
$ws.Range("H20").Value = 20.17001148208421
$ws.Range("J20").Value = 11.44199572254184
$ws.Range("K20").Value = 11.92731385870579
$ws.Range("L20").Value = 9.639252325193681
$ws.Range("B21").Value = 17.8308707945248
$ws.Range("D21").Value = 7.732023423101247
$ws.Range("E21").Value = 14.40696109637167
$ws.Range("F21").Value = 44.14483872517769
$ws.Range("G21").Value = 53.50609247338555
$ws.Range("H21").Value = 20.15836559301025
$ws.Range("J21").Value = 11.42282628287319
$ws.Range("K21").Value = 12.16764517706277
$ws.Range("L21").Value = 9.675148507066938
$ws.Range("B22").Value = 17.91629191845121
$ws.Range("D22").Value = 7.729276939826876
$ws.Range("E22").Value = 14.39220171605458
$ws.Range("F22").Value = 44.19951686755271
$ws.Range("G22").Value = 53.66499493408564
$ws.Range("H22").Value = 20.15355372839413
$ws.Range("J22").Value = 11.41080252328886
$ws.Range("K22").Value = 12.32441717527018
$ws.Range("L22").Value = 9.699952408616788
$ws.Range("B23").Value = 17.87047516285972
$ws.Range("D23").Value = 7.730694198034723
$ws.Range("E23").Value = 14.39999862388141
$ws.Range("F23").Value = 44.16965551848796
$ws.Range("G23").Value = 53.57926529643957
$ws.Range("H23").Value = 20.15586577035966
$ws.Range("J23").Value = 11.41717377261794
$ws.Range("K23").Value = 12.24078993925983
$ws.Range("L23").Value = 9.686592761664055
$ws.Range("B24").Value = 17.70202883932507
$ws.Range("D24").Value = 7.7371318678195
$ws.Range("E24").Value = 14.43129804256451
$ws.Range("F24").Value = 44.07150577882516
$ws.Range("G24").Value = 53.27494366856178
$ws.Range("H24").Value = 20.17024908344303
$ws.Range("J24").Value = 11.44231758039123
$ws.Range("K24").Value = 11.92338472150143
$ws.Range("L24").Value = 9.638688656033109
$ws.Range("B25").Value = 17.53188827386623
$ws.Range("D25").Value = 7.746377805510979
$ws.Range("E25").Value = 14.46885628728553
$ws.Range("F25").Value = 43.9977685126988
$ws.Range("G25").Value = 52.99104112629395
$ws.Range("H25").Value = 20.19777787098662
$ws.Range("J25").Value = 11.47161698936926
$ws.Range("K25").Value = 11.58156127104012
$ws.Range("L25").Value = 9.592939607993239
